$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a plain-text value into a cell without Excel's automatic
# "looks like a number" coercion (which would turn e.g. "249.30" into the
# number 249.3 and lose the trailing zero / text semantics). We do this by
# writing a text formula that yields the literal string, then converting
# that formula to a static value in place via copy / paste-special-values.
# This keeps the cell's style/number-format untouched (no quote-prefix or
# "@" format side effects), matching a plain inline-string text cell.
function Set-TextValue {
    param(
        [string]$CellRef,
        [string]$Value
    )
    $escaped = $Value -replace '"', '""'
    $cell = $ws.Range($CellRef)
    $cell.Formula = '="' + $escaped + '"'
    $cell.Copy() | Out-Null
    $cell.PasteSpecial(-4163) | Out-Null  # xlPasteValues
}

Set-TextValue "D2"  "249.30"
Set-TextValue "D3"  "22.00"
Set-TextValue "D4"  "5.354"
Set-TextValue "D5"  "0.05622"
Set-TextValue "D6"  "3.431"
Set-TextValue "D7"  "6.382"
Set-TextValue "D8"  "0.8171"
Set-TextValue "D9"  "0.9288"
Set-TextValue "D10" "0.1450"
Set-TextValue "D11" "0.07481"
Set-TextValue "D12" "0.03239"
Set-TextValue "D14" "0.09317"
Set-TextValue "D15" "3.554"
Set-TextValue "D16" "0.001636"
Set-TextValue "D17" "0.04729"
Set-TextValue "D18" "0.0005759"
Set-TextValue "D19" "0.006385"
Set-TextValue "D20" "0.005067"
Set-TextValue "D23" "3.735"
Set-TextValue "D24" "2.162"
Set-TextValue "D25" "0.3305"
Set-TextValue "D26" "0.1328"
Set-TextValue "E27" "26AAXTokenAABWorstin24h"
Set-TextValue "D40" "0.03944"
Set-TextValue "D41" "0.006829"
Set-TextValue "E41" "40KickTokenKICK"
Set-TextValue "D42" "0.1067"
Set-TextValue "D43" "0.003400"
Set-TextValue "D44" "0.008559"
Set-TextValue "D45" "0.00005571"
Set-TextValue "D49" "0.1954"
Set-TextValue "D50" "0.00002100"

$excel.CutCopyMode = $false
